$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Training Dashboard"
#   For every data row (3..31):
#     - column H (PERIOD TO EXPIRE) decreases by 1
#     - column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025
# ---------------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item(1)

for ($row = 3; $row -le 31; $row++) {
    $hCell = $wsTraining.Cells.Item($row, 8)   # column H
    $iCell = $wsTraining.Cells.Item($row, 9)   # column I

    $currentPeriod = $hCell.Value()
    $hCell.Value = $currentPeriod - 1

    # Force the date to be written as literal text (matches the workbook's
    # existing convention of storing these dates as plain strings, not as
    # Excel date serials).
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
}

# ---------------------------------------------------------------------------
# Sheet 2: "Exam Dashboard"
#   - widen column B from 17 to 20 characters
#   - insert a new exam result row (row 5) before the TOTAL AVERAGE row
#   - update the TOTAL AVERAGE value (now on row 6)
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item(2)

# Widen column B. Excel snaps ColumnWidth to the workbook's default font
# character grid, so 19.1 is the value that round-trips to a stored
# width of exactly 20.
$wsExam.Columns.Item(2).ColumnWidth = 19.1

# Push the existing TOTAL AVERAGE row (row 5) down to row 6 and create a
# fresh blank row 5 for the new exam entry.
$wsExam.Rows.Item(5).Insert()

$wsExam.Range("A5").Value = 3
$wsExam.Range("B5").Value = "Asrs Stacker Crane"

$wsExam.Range("C5").NumberFormat = "@"
$wsExam.Range("C5").Value = "30-Oct-2025"

$wsExam.Range("D5").NumberFormat = "@"
$wsExam.Range("D5").Value = "83.54%"

$wsExam.Range("E5").Value = "VALID"
$wsExam.Range("F5").Value = "Approved Score. date is valid"

# Update the TOTAL AVERAGE (now shifted to row 6)
$wsExam.Range("D6").NumberFormat = "@"
$wsExam.Range("D6").Value = "87.02%"
